# Apply the "Add files via upload" revision to the ESG Topics worksheet.
#
# Summary of the change:
#  - A1 header text loses its trailing space: "Catégorie " -> "Catégorie"
#  - Two new header columns are inserted right after "Thème" (column B):
#       C1 = "Impact produit (risques)"
#       D1 = "Demande client"
#    and the previous C1:F1 headers ("Évolution du risque", "Demande",
#    "Besoin d'innovation", "Pression externe") shift right by two columns.
#    NOTE: only the header row shifts - the data rows (2-13) keep using
#    columns C:H, only their numeric contents change.
#  - The two header strings that fall out of use ("Évolution du risque",
#    "Demande") are no longer referenced by any cell.
#  - Every data row's scores (columns C:G) are updated to new values; the
#    H column keeps its existing SUM(C:G) formulas and recalculates.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -------------------------------------------------------
# Write right-to-left so we never clobber a value before it has been
# copied to its new home column.
$ws.Range("H1").Value = "Score total"
$ws.Range("G1").Value = "Maturité VA"
$ws.Range("F1").Value = "Pression externe"
$ws.Range("E1").Value = "Besoin d’innovation"
$ws.Range("D1").Value = "Demande client"
$ws.Range("C1").Value = "Impact produit (risques)"
$ws.Range("B1").Value = "Thème"
$ws.Range("A1").Value = "Catégorie"

# --- Data rows: updated scores (columns C:G); H keeps its SUM formula --
$ws.Range("C2").Value = 5
$ws.Range("D2").Value = 4
$ws.Range("E2").Value = 4
$ws.Range("F2").Value = 3
$ws.Range("G2").Value = 5

$ws.Range("C3").Value = 3
$ws.Range("D3").Value = 2
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 2
$ws.Range("G3").Value = 5

$ws.Range("C4").Value = 4
$ws.Range("D4").Value = 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 4
$ws.Range("G4").Value = 3

$ws.Range("C5").Value = 3
$ws.Range("D5").Value = 4
$ws.Range("E5").Value = 5
$ws.Range("F5").Value = 2
$ws.Range("G5").Value = 5

$ws.Range("C6").Value = 3
$ws.Range("D6").Value = 2
$ws.Range("E6").Value = 4
$ws.Range("F6").Value = 3
$ws.Range("G6").Value = 4

$ws.Range("C7").Value = 3
$ws.Range("D7").Value = 4
$ws.Range("E7").Value = 4
$ws.Range("F7").Value = 2
$ws.Range("G7").Value = 3

$ws.Range("C8").Value = 2
$ws.Range("D8").Value = 5
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 4
$ws.Range("G8").Value = 2

$ws.Range("C9").Value = 4
$ws.Range("D9").Value = 3
$ws.Range("E9").Value = 5
$ws.Range("F9").Value = 4
$ws.Range("G9").Value = 5

$ws.Range("C10").Value = 3
$ws.Range("D10").Value = 3
$ws.Range("E10").Value = 5
$ws.Range("F10").Value = 3
$ws.Range("G10").Value = 5

$ws.Range("C11").Value = 3
$ws.Range("D11").Value = 2
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 3
$ws.Range("G11").Value = 1

$ws.Range("C12").Value = 5
$ws.Range("D12").Value = 4
$ws.Range("E12").Value = 5
$ws.Range("F12").Value = 2
$ws.Range("G12").Value = 4

$ws.Range("C13").Value = 1
$ws.Range("D13").Value = 1
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 2
$ws.Range("G13").Value = 2

# --- Match the saved selection shown in the diff -----------------------
$ws.Range("C8").Select()
